# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.709.57"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.598.76"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'211.21"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.823.08"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "1.597.33"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "'65.30"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "26.681.62"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").Value = "'210.05"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "'7.15"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'2.30"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'143.24"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'15.32"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'0.0519"
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "'2.97"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").Value = "1.289.31"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "'0.619"
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'0.0172"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  +17.08%  "
$ws.Range("D40").Value = "'0.826"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'63.09"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "1.740.13"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").Value = "'91.26"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  -2.06%  "
